$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 'Muhammad Luqman_20251202_121836'
$ws.Range("C12").Value = 'Muhammad Luqman'
$ws.Range("D12").Value = 18
$ws.Range("E12").Value = 'Male'
$ws.Range("F12").Value = '2025-12-02 12:18:37'
$ws.Range("H12").Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range("J12").Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range("K12").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("M12").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N12").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("P12").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("Q12").Value = 'Velveeta Original Shells & Cheese (microwave cups)'
$ws.Range("S12").Value = 'Muy cremoso, porción individual, rápido, salado, ideal para niños'
$ws.Range("T12").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("V12").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("W12").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("Y12").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("Z12").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AB12").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AC12").Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AE12").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'
$ws.Range("AF12").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AH12").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = ''
$ws.Range("B12").ClearFormats()
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = '0.579'
$ws.Range("I12").ClearFormats()
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = '0.567'
$ws.Range("L12").ClearFormats()
$ws.Range("O12").NumberFormat = "@"
$ws.Range("O12").Value = '0.556'
$ws.Range("O12").ClearFormats()
$ws.Range("R12").NumberFormat = "@"
$ws.Range("R12").Value = '0.600'
$ws.Range("R12").ClearFormats()
$ws.Range("U12").NumberFormat = "@"
$ws.Range("U12").Value = '0.592'
$ws.Range("U12").ClearFormats()
$ws.Range("X12").NumberFormat = "@"
$ws.Range("X12").Value = '0.526'
$ws.Range("X12").ClearFormats()
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = '0.660'
$ws.Range("AA12").ClearFormats()
$ws.Range("AD12").NumberFormat = "@"
$ws.Range("AD12").Value = '0.659'
$ws.Range("AD12").ClearFormats()
$ws.Range("AG12").NumberFormat = "@"
$ws.Range("AG12").Value = '0.647'
$ws.Range("AG12").ClearFormats()

$ws.Range("G12").Value = '{
  "portion": 0.8,
  "diet": 0.42857142857142855,
  "salt": 0.8,
  "fat": 0.6,
  "natural": 0.4,
  "convenience": 0.8,
  "price": 1.0
}'

$ws.Rows(12).AutoFit()
